$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd34"
$ws.Range("C2").Value = "Selp"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 106.3055973333333
$ws.Range("H2").Value = 318.916792
$ws.Range("I2").Value = 0.2547398208373942
$ws.Range("J2").Value = 0.2547398208373942
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 17.47090666666666
$ws.Range("N2").Value = 52.41271999999999
$ws.Range("O2").Value = 0.9803569739482672
$ws.Range("P2").Value = 0.9803569739482673
$ws.Range("Q2").Value = 1857.255169154915
$ws.Range("R2").Value = 16715.29652239424
$ws.Range("S2").Value = 0.2497359599002715
$ws.Range("T2").Value = 0.2497359599002715

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd34"
$ws.Range("C3").Value = "Selp"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 106.3055973333333
$ws.Range("H3").Value = 318.916792
$ws.Range("I3").Value = 0.2547398208373942
$ws.Range("J3").Value = 0.2547398208373942
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.16459
$ws.Range("N3").Value = 0.49377
$ws.Range("O3").Value = 0.0092357516081294
$ws.Range("P3").Value = 0.0092357516081294
$ws.Range("Q3").Value = 17.49683826509333
$ws.Range("R3").Value = 157.47154438584
$ws.Range("S3").Value = 0.002352713709953558
$ws.Range("T3").Value = 0.002352713709953558

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cd34"
$ws.Range("C4").Value = "Selp"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 106.3055973333333
$ws.Range("H4").Value = 318.916792
$ws.Range("I4").Value = 0.2547398208373942
$ws.Range("J4").Value = 0.2547398208373942
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1854676666666667
$ws.Range("N4").Value = 0.556403
$ws.Range("O4").Value = 0.01040727444360334
$ws.Range("P4").Value = 0.01040727444360334
$ws.Range("Q4").Value = 19.71625109101955
$ws.Range("R4").Value = 177.446259819176
$ws.Range("S4").Value = 0.002651147227169106
$ws.Range("T4").Value = 0.002651147227169106

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cd34"
$ws.Range("C5").Value = "Selp"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 296.1091513333333
$ws.Range("H5").Value = 888.327454
$ws.Range("I5").Value = 0.7095655736964096
$ws.Range("J5").Value = 0.7095655736964096
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 17.47090666666666
$ws.Range("N5").Value = 52.41271999999999
$ws.Range("O5").Value = 0.9803569739482672
$ws.Range("P5").Value = 0.9803569739482673
$ws.Range("Q5").Value = 5173.295346090542
$ws.Range("R5").Value = 46559.65811481488
$ws.Range("S5").Value = 0.6956275586468783
$ws.Range("T5").Value = 0.6956275586468783

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cd34"
$ws.Range("C6").Value = "Selp"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 296.1091513333333
$ws.Range("H6").Value = 888.327454
$ws.Range("I6").Value = 0.7095655736964096
$ws.Range("J6").Value = 0.7095655736964096
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.16459
$ws.Range("N6").Value = 0.49377
$ws.Range("O6").Value = 0.0092357516081294
$ws.Range("P6").Value = 0.0092357516081294
$ws.Range("Q6").Value = 48.73660521795333
$ws.Range("R6").Value = 438.62944696158
$ws.Range("S6").Value = 0.006553371388339875
$ws.Range("T6").Value = 0.006553371388339875

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cd34"
$ws.Range("C7").Value = "Selp"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 296.1091513333333
$ws.Range("H7").Value = 888.327454
$ws.Range("I7").Value = 0.7095655736964096
$ws.Range("J7").Value = 0.7095655736964096
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1854676666666667
$ws.Range("N7").Value = 0.556403
$ws.Range("O7").Value = 0.01040727444360334
$ws.Range("P7").Value = 0.01040727444360334
$ws.Range("Q7").Value = 54.91867337644022
$ws.Range("R7").Value = 494.268060387962
$ws.Range("S7").Value = 0.007384643661191388
$ws.Range("T7").Value = 0.007384643661191388

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cd34"
$ws.Range("C8").Value = "Selp"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.895733
$ws.Range("H8").Value = 44.68719899999999
$ws.Range("I8").Value = 0.03569460546619627
$ws.Range("J8").Value = 0.03569460546619627
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 17.47090666666666
$ws.Range("N8").Value = 52.41271999999999
$ws.Range("O8").Value = 0.9803569739482672
$ws.Range("P8").Value = 0.9803569739482673
$ws.Range("Q8").Value = 260.2419609745866
$ws.Range("R8").Value = 2342.177648771279
$ws.Range("S8").Value = 0.03499345540111745
$ws.Range("T8").Value = 0.03499345540111745

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cd34"
$ws.Range("C9").Value = "Selp"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.895733
$ws.Range("H9").Value = 44.68719899999999
$ws.Range("I9").Value = 0.03569460546619627
$ws.Range("J9").Value = 0.03569460546619627
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.16459
$ws.Range("N9").Value = 0.49377
$ws.Range("O9").Value = 0.0092357516081294
$ws.Range("P9").Value = 0.0092357516081294
$ws.Range("Q9").Value = 2.45168869447
$ws.Range("R9").Value = 22.06519825023
$ws.Range("S9").Value = 0.0003296665098359667
$ws.Range("T9").Value = 0.0003296665098359666

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cd34"
$ws.Range("C10").Value = "Selp"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 14.895733
$ws.Range("H10").Value = 44.68719899999999
$ws.Range("I10").Value = 0.03569460546619627
$ws.Range("J10").Value = 0.03569460546619627
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1854676666666667
$ws.Range("N10").Value = 0.556403
$ws.Range("O10").Value = 0.01040727444360334
$ws.Range("P10").Value = 0.01040727444360334
$ws.Range("Q10").Value = 2.762676842799666
$ws.Range("R10").Value = 24.864091585197
$ws.Range("S10").Value = 0.0003714835552428486
$ws.Range("T10").Value = 0.0003714835552428486
